# Updates cryptos list values per the Jun 17 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.495.69"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "1.737.91"
$ws.Range("E3").Value = "  +4.44%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'245.35"
$ws.Range("E5").Value = "  +4.93%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4790"
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("D8").Value = "'0.2679"
$ws.Range("E8").Value = "  +4.17%  "
$ws.Range("D9").Value = "'0.06240"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "1.739.54"
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("D11").Value = "'0.07131"
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").Value = "'15.79"
$ws.Range("E12").Value = "  +8.24%  "
$ws.Range("D13").Value = "'0.6168"
$ws.Range("E13").Value = "  +8.74%  "
$ws.Range("D14").Value = "'4.550"
$ws.Range("E14").Value = "  +4.49%  "
$ws.Range("D15").Value = "'77.10"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "26.503.45"
$ws.Range("E17").Value = "  +4.03%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "'0.000006901"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").Value = "1.963.36"
$ws.Range("D22").Value = "'4.600"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").Value = "'8.869"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "'5.354"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").Value = "'136.15"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +3.35%  "
$ws.Range("D27").Value = "'1.805"
$ws.Range("E27").Value = "  +6.63%  "
$ws.Range("E28").Value = "  +4.59%  "
$ws.Range("D29").Value = "'107.09"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").Value = "'3.985"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "'3.731"
$ws.Range("E31").Value = "  +3.51%  "
$ws.Range("D32").Value = "'0.07856"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").Value = "'0.04563"
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").Value = "'2.621"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6375"
$ws.Range("E35").Value = "  +6.95%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9981"
$ws.Range("E36").Value = "  +6.19%  "
$ws.Range("D37").Value = "'0.9290"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").Value = "'112.84"
$ws.Range("E38").Value = "  +9.96%  "
$ws.Range("D39").Value = "'2.423"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").Value = "'1.975"
$ws.Range("E40").Value = "  +9.31%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "'5.768"
$ws.Range("E42").Value = "  +17.58%  "
$ws.Range("D43").Value = "'0.01504"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("D44").Value = "'0.3899"
$ws.Range("E44").Value = "  +5.60%  "
$ws.Range("D45").Value = "'6.844"
$ws.Range("E45").Value = "  +12.06%  "
$ws.Range("D46").Value = "'0.1206"
$ws.Range("E46").Value = "  +9.42%  "
$ws.Range("D47").Value = "'0.05331"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "'30.74"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("D49").Value = "'7.849"
$ws.Range("E49").Value = "  +6.58%  "
$ws.Range("D50").Value = "'1.256"
$ws.Range("E50").Value = "  +5.36%  "
$ws.Range("D51").Value = "'0.3436"
$ws.Range("E51").Value = "  +4.65%  "
